# Update the "regula falsa" (false position) table on Hoja1 with a new
# iteration run (f(x) = x^2 - 2 over [-2, 1]) that converges in 7 steps
# instead of 10, so rows 9-11 are no longer needed and the used range
# shrinks from A1:H11 to A1:H8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-8 (columns A:H -> n, x_m, x_i, x_s, f_m, f_i, f_s, E)
$data = @(
    @(1, 0,               -2, 1,               -2,               2, -1,               1.0005),
    @(2, -1,               -2, 0,               -1,               2, -2,               1),
    @(3, -1.33333333333333,-2, -1,               -0.222222222222222, 2, -1,             0.333333333333333),
    @(4, -1.4,             -2, -1.33333333333333,-0.0400000000000003,2, -0.222222222222222,0.0666666666666664),
    @(5, -1.41176470588235,-2, -1.4,             -0.0069204152249138,2, -0.0400000000000003,0.0117647058823529),
    @(6, -1.41379310344828,-2, -1.41176470588235,-0.001189060642093, 2, -0.0069204152249138,0.0020283975659229),
    @(7, -1.41414141414141,-2, -1.41379310344828,-0.0002040608101208,2, -0.001189060642093,0.0003483106931385)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($row, $c).Value = $values[$c - 1]
    }
}

# Remove the now-unused rows 9-11 (table shrank from 10 to 7 iterations)
$ws.Range("A9:H11").Delete()
